$d = $word.ActiveDocument

# Fix the Functor definition: Applicative should require Pointed as well as Functor.
$d.Content.Find.Execute(
    "category Applicative<p; Functor<p>> =",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "category Applicative<p; Pointed<p>; Functor<p>> =",
    2)
